$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (slash -> hyphen format)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

# Rows whose day-number is <= 12, so "dd-mm-yyyy" reads as an ambiguous,
# parseable date (e.g. "01-08-2022") - Excel's Value setter would silently
# convert those to a date serial unless we pin the cell to Text format
# first. Unambiguous ones (day > 12) are safe to assign directly, keeping
# their style untouched.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    if ($ambiguousRows -contains $row) {
        $cell.NumberFormat = "@"
        $cell.Value = $dates[$row]
        # Restore the default "Normal" style now that the literal text is
        # committed, so the cell doesn't end up tagged with a lingering
        # Text number format that the source workbook never had.
        $cell.Style = "Normal"
    }
    else {
        $cell.Value = $dates[$row]
    }
}

# Row 3 additionally gets its attendance counts updated
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1
